# Update MSME country indicator figures on the "Summary" sheet.
#
# The source values are stored as text (not numbers) in the workbook, so we
# force each target cell to a text number-format before writing the new
# string - otherwise Excel's normal type-inference would turn a value like
# "2.73" into a numeric cell, which does not match the original data model
# (every one of these cells is backed by a shared string, t="s").
#
# Cell map (row 9 header is Micro / SMEs / MSMEs -> columns B / C / D):
#   Row 11 "Enterprises density (per 1000 people)": B 2.7 -> 2.73, C 1 -> 0.97
#   Row 12 "Enterprises (% of total)":              B 71.5 -> 71.54, C 25.4 -> 25.39, D 96.9 -> 96.93
#   Row 16 "Value added to the economy (% of total)": B 19.9 -> 19.93, D 44.3 -> 44.33
#
# (The duplicate "2.73" shared-string entry used elsewhere in the sheet -
# e.g. the "Postage and communications" row - is naturally deduplicated by
# the engine once B11 becomes "2.73" too, which is exactly what collapses
# the shared-string table from 362 to 361 unique entries.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2.73"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "0.97"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "71.54"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "25.39"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.93"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "19.93"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "44.33"
